$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    "B2" = 8.900953978797508
    "C2" = 5.853951952849987
    "E2" = 24.40626903363801
    "F2" = 39.50262735724699
    "G2" = 24.30270867600074
    "H2" = 13.29341888171801
    "J2" = 7.727426390497756
    "K2" = 8.091957464186308
    "N2" = 17.51430069998071
    "O2" = 19.61826106759793
    "B3" = 8.588643214191647
    "C3" = 5.74537355482572
    "E3" = 24.08341760861649
    "F3" = 39.31606069416119
    "G3" = 24.41093236236594
    "H3" = 13.34188378322304
    "J3" = 7.747786845177074
    "K3" = 7.872516146082173
    "N3" = 17.56232177769927
    "O3" = 19.70299983665339
    "B4" = 8.392033943782351
    "C4" = 5.677273461581621
    "E4" = 23.88901999585932
    "F4" = 39.2120571924717
    "G4" = 24.48545411430007
    "H4" = 13.37358725792813
    "J4" = 7.76101686104019
    "K4" = 7.735407093935599
    "N4" = 17.59340333480534
    "O4" = 19.75903571811734
    "B5" = 8.310817553697584
    "C5" = 5.649186543081151
    "E5" = 23.81085916733012
    "F5" = 39.17236032852504
    "G5" = 24.51784128644061
    "H5" = 13.38699632658563
    "J5" = 7.766591840888973
    "K5" = 7.679017748233003
    "N5" = 17.60647170255538
    "O5" = 19.78287685587882
    "B6" = 8.297269223958947
    "C6" = 5.644503202653018
    "E6" = 23.79794720173647
    "F6" = 39.16593176194472
    "G6" = 24.52334076970858
    "H6" = 13.3892524761725
    "J6" = 7.767528666049806
    "K6" = 7.66962558869798
    "N6" = 17.60866602591593
    "O6" = 19.78689638574078
    "B7" = 8.390942900079594
    "C7" = 5.676895997043825
    "E7" = 23.88796148664375
    "F7" = 39.21151091346778
    "G7" = 24.48588273889524
    "H7" = 13.37376611429321
    "J7" = 7.761091303045872
    "K7" = 7.73464858875664
    "N7" = 17.59357794889798
    "O7" = 19.75935317620304
    "B8" = 8.794345820752527
    "C8" = 5.816825913542958
    "E8" = 24.29421073728707
    "F8" = 39.43613056347474
    "G8" = 24.33834245762995
    "H8" = 13.30972593274272
    "J8" = 7.73429570718434
    "K8" = 8.016832694350915
    "N8" = 17.53052760059903
    "O8" = 19.64664680157586
    "B9" = 9.542125317731655
    "C9" = 6.078832218979755
    "E9" = 25.11709845647103
    "F9" = 39.95868411586758
    "G9" = 24.1135331110803
    "H9" = 13.19956671462303
    "J9" = 7.687511541431689
    "K9" = 8.552176603285346
    "N9" = 17.4195096944185
    "O9" = 19.45747196327265
    "B10" = 10.0594618736386
    "C10" = 6.262472415175492
    "E10" = 25.7320010169767
    "F10" = 40.39019867366212
    "G10" = 23.98831624404232
    "H10" = 13.12801317167199
    "J10" = 7.65662491990407
    "K10" = 9.023630319276384
    "N10" = 17.34558000803041
    "O10" = 19.33797496892892
    "B11" = 10.28686439652517
    "C11" = 6.343840679907096
    "E11" = 26.01278023280308
    "F11" = 40.59624665979749
    "G11" = 23.9401512489428
    "H11" = 13.09749411336063
    "J11" = 7.64332501029374
    "K11" = 9.228663117855644
    "N11" = 17.31359267448851
    "O11" = 19.28786197727724
    "B12" = 10.37176524537266
    "C12" = 6.37432120033717
    "E12" = 26.1191575854815
    "F12" = 40.67561891015396
    "G12" = 23.92318693359406
    "H12" = 13.08622913820334
    "J12" = 7.638396189562163
    "K12" = 9.304912842598529
    "N12" = 17.30171533184856
    "O12" = 19.26949772330111
    "B13" = 10.35353523940168
    "C13" = 6.367771736176008
    "E13" = 26.09624655324051
    "F13" = 40.65846571183356
    "G13" = 23.92678365889813
    "H13" = 13.0886422703209
    "J13" = 7.639452920999658
    "K13" = 9.288553504067428
    "N13" = 17.30426286560946
    "O13" = 19.27342552685004
    "B14" = 10.29387389339756
    "C14" = 6.346355090577868
    "E14" = 26.02153147511011
    "F14" = 40.60275000829289
    "G14" = 23.93872998310038
    "H14" = 13.09656148702574
    "J14" = 7.642917359524414
    "K14" = 9.234964304497364
    "N14" = 17.31261080172138
    "O14" = 19.28633885449279
    "B15" = 10.25716987961757
    "C15" = 6.333192992202084
    "E15" = 25.9757702313708
    "F15" = 40.56879616670054
    "G15" = 23.94621373901421
    "H15" = 13.10145025036117
    "J15" = 7.64505342574451
    "K15" = 9.201957170045077
    "N15" = 17.3177548081907
    "O15" = 19.29432845211766
    "B16" = 10.04443467080162
    "C16" = 6.257109413109253
    "E16" = 25.71366413849735
    "F16" = 40.37692445025443
    "G16" = 23.99164179904859
    "H16" = 13.13004852173848
    "J16" = 7.657509172364425
    "K16" = 9.010037984846111
    "N16" = 17.34770346218407
    "O16" = 19.34133557749243
    "B17" = 9.911842941063103
    "C17" = 6.209864212385603
    "E17" = 25.55306971430068
    "F17" = 40.2616758289915
    "G17" = 24.02177061060587
    "H17" = 13.14811272888933
    "J17" = 7.665342347217003
    "K17" = 8.889858844317517
    "N17" = 17.36649640366372
    "O17" = 19.37126188071993
    "B18" = 9.834836279977441
    "C18" = 6.182486818081689
    "E18" = 25.46080538175089
    "F18" = 40.19630909214074
    "G18" = 24.0399273484459
    "H18" = 13.15869395747129
    "J18" = 7.669918450819889
    "K18" = 8.819848770903878
    "N18" = 17.37746036092708
    "O18" = 19.38887431991436
    "B19" = 9.808637754660172
    "C19" = 6.173183024606701
    "E19" = 25.42958742491546
    "F19" = 40.17433690909499
    "G19" = 24.04621671647059
    "H19" = 13.16230941820988
    "J19" = 7.671479989470277
    "K19" = 8.795993513059836
    "N19" = 17.3811991719899
    "O19" = 19.39490617108336
    "B20" = 9.926035033814506
    "C20" = 6.214914722612589
    "E20" = 25.57015509185068
    "F20" = 40.27384929281998
    "G20" = 24.01847763770809
    "H20" = 13.14616997833113
    "J20" = 7.664501181103425
    "K20" = 8.902744109627466
    "N20" = 17.3644798522215
    "O20" = 19.36803480263118
    "B21" = 10.31143127506807
    "C21" = 6.352654836671497
    "E21" = 26.04347646673493
    "F21" = 40.61907897859082
    "G21" = 23.9351863872689
    "H21" = 13.09422749914004
    "J21" = 7.64189685350865
    "K21" = 9.25074277789381
    "N21" = 17.31015242353084
    "O21" = 19.28252926130217
    "B22" = 10.55622154184769
    "C22" = 6.440733268180248
    "E22" = 26.35307998851953
    "F22" = 40.85252641953199
    "G22" = 23.88818495851163
    "H22" = 13.06198163115472
    "J22" = 7.627750432569148
    "K22" = 9.470055049951215
    "N22" = 17.27601901474286
    "O22" = 19.23021709762675
    "B23" = 10.42624196736941
    "C23" = 6.393908106796846
    "E23" = 26.18784723071284
    "F23" = 40.72723468542917
    "G23" = 23.91258718492166
    "H23" = 13.07903620341376
    "J23" = 7.635243409225184
    "K23" = 9.35375771796555
    "N23" = 17.29411131110691
    "O23" = 19.25780978093215
    "B24" = 9.919621206255739
    "C24" = 6.212632056683152
    "E24" = 25.56243058707557
    "F24" = 40.26834288942241
    "G24" = 24.01996378962572
    "H24" = 13.14704768569731
    "J24" = 7.664881245799338
    "K24" = 8.896921536253533
    "N24" = 17.36539103777525
    "O24" = 19.36949249569494
    "B25" = 9.345092678881093
    "C25" = 6.009418274307518
    "E25" = 24.89226177160728
    "F25" = 39.95868411586758
    "G25" = 24.1135331110803
    "H25" = 13.22771887200894
    "J25" = 7.699553846725362
    "K25" = 8.40719023966923
    "N25" = 17.44819785015589
    "O25" = 19.50523189779474
}

foreach ($cell in $data.Keys) {
    $ws.Range($cell).Value = $data[$cell]
}
